# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the a00a7228-422b-48f1-b114-67c1f80c027f
# entry across the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
# Row 4 corresponds to a00a7228-422b-48f1-b114-67c1f80c027f.md
# Column G = "Latest HO Xliff Generate Date"
$overview.Range("G4").Value = "2016-10-19 16:50:55"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
# Row 4 corresponds to a00a7228-422b-48f1-b114-67c1f80c027f.*.zh-cn.xlf
# Column H = "Correspond Handoff Datetime"
$zhcn.Range("H4").Value = "2016-10-19 16:50:44"
# Column K = "Correspond Handback DateTime"
$zhcn.Range("K4").Value = "2016-10-19 16:51:30"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
# Row 4 corresponds to a00a7228-422b-48f1-b114-67c1f80c027f.*.de-de.xlf
# Column H = "Correspond Handoff Datetime" (shares same value/shared string as Overview!G4)
$dede.Range("H4").Value = "2016-10-19 16:50:55"
# Column K = "Correspond Handback DateTime"
$dede.Range("K4").Value = "2016-10-19 16:51:48"
